$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.660.53'
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.503.24'
$ws.Range("E3").Value = '  -0.04%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.21'
$ws.Range("E5").Value = '  -0.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.86'
$ws.Range("E6").Value = '  +0.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.501.86'
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("E10").Value = '  +1.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.167'
$ws.Range("E11").Value = '  +0.18%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.358'
$ws.Range("E12").Value = '  +3.38%  '
$ws.Range("E13").Value = '  +1.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.961.17'
$ws.Range("E14").Value = '  -0.04%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '69.545.39'
$ws.Range("E15").Value = '  +0.26%  '
$ws.Range("E16").Value = '  +1.81%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '24.75'
$ws.Range("E17").Value = '  -0.73%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.500.20'
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.24'
$ws.Range("E19").Value = '  -0.81%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.50'
$ws.Range("E20").Value = '  -4.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '348.42'
$ws.Range("E21").Value = '  -0.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.92'
$ws.Range("E22").Value = '  -0.74%  '
$ws.Range("E23").Value = '  -0.31%  '
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.55'
$ws.Range("E25").Value = '  +2.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.96'
$ws.Range("E26").Value = '  -0.55%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.78'
$ws.Range("E27").Value = '  -1.57%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.630.26'
$ws.Range("E28").Value = '  -0.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.27%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0893'
$ws.Range("E30").Value = '  -1.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.83'
$ws.Range("E31").Value = '  -0.93%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '458.10'
$ws.Range("E32").Value = '  -0.94%  '
$ws.Range("E33").Value = '  -3.32%  '
$ws.Range("E34").Value = '  -0.96%  '
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("E36").Value = '  -1.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '156.91'
$ws.Range("E37").Value = '  +1.70%  '
$ws.Range("E38").Value = '  +0.60%  '
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("E41").Value = '  -0.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.69'
$ws.Range("E42").Value = '  -1.14%  '
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("E44").Value = '  +0.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.21'
$ws.Range("E45").Value = '  -3.87%  '
$ws.Range("E46").Value = '  -6.35%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '141.89'
$ws.Range("E47").Value = '  -0.93%  '
$ws.Range("E48").Value = '  -0.41%  '
$ws.Range("E49").Value = '  -1.77%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0732'
$ws.Range("E50").Value = '  +0.22%  '
$ws.Range("E51").Value = '  -0.53%  '
